$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before row 403; this pushes the existing
# rows 403-419 down to 405-421, preserving their data and formatting.
$ws.Rows.Item(403).Resize(2).Insert()

# Row 403 (new): Femacal de La Calera - Coliflor - Primera
$ws.Cells.Item(403, 1).Value = 3
$ws.Cells.Item(403, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(403, 3).Value = "Coquimbo"
$ws.Cells.Item(403, 4).Value = 44509
$ws.Cells.Item(403, 5).Value = 5
$ws.Cells.Item(403, 6).Value = 100112008
$ws.Cells.Item(403, 7).Value = "Coliflor"
$ws.Cells.Item(403, 8).Value = "Sin especificar"
$ws.Cells.Item(403, 9).Value = "Primera"
$ws.Cells.Item(403, 10).Value = 2550
$ws.Cells.Item(403, 11).Value = 600
$ws.Cells.Item(403, 12).Value = 650
$ws.Cells.Item(403, 13).Value = 625
$ws.Cells.Item(403, 14).Value = "`$/unidad"
$ws.Cells.Item(403, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(403, 16).Value = 625
$ws.Cells.Item(403, 17).Value = 1
$ws.Cells.Item(403, 18).Value = "Hortaliza"

# Row 404 (new): Femacal de La Calera - Coliflor - Segunda
$ws.Cells.Item(404, 1).Value = 3
$ws.Cells.Item(404, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(404, 3).Value = "Coquimbo"
$ws.Cells.Item(404, 4).Value = 44509
$ws.Cells.Item(404, 5).Value = 5
$ws.Cells.Item(404, 6).Value = 100112008
$ws.Cells.Item(404, 7).Value = "Coliflor"
$ws.Cells.Item(404, 8).Value = "Sin especificar"
$ws.Cells.Item(404, 9).Value = "Segunda"
$ws.Cells.Item(404, 10).Value = 1200
$ws.Cells.Item(404, 11).Value = 500
$ws.Cells.Item(404, 12).Value = 500
$ws.Cells.Item(404, 13).Value = 500
$ws.Cells.Item(404, 14).Value = "`$/unidad"
$ws.Cells.Item(404, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(404, 16).Value = 500
$ws.Cells.Item(404, 17).Value = 1
$ws.Cells.Item(404, 18).Value = "Hortaliza"
